$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 23.301814
$ws.Range("H2").Value = 46.603628
$ws.Range("I2").Value = 0.1534311298773869
$ws.Range("J2").Value = 0.1083150512651137
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.704179666666667
$ws.Range("N2").Value = 5.112539
$ws.Range("O2").Value = 0.3144996488703566
$ws.Range("P2").Value = 0.3144996488703566
$ws.Range("Q2").Value = 39.71047761524866
$ws.Range("R2").Value = 238.262865691492
$ws.Range("S2").Value = 0.04825403647222025
$ws.Range("T2").Value = 0.03406504559025293

$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 23.301814
$ws.Range("H3").Value = 46.603628
$ws.Range("I3").Value = 0.1534311298773869
$ws.Range("J3").Value = 0.1083150512651137
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.714521666666666
$ws.Range("N3").Value = 11.143565
$ws.Range("O3").Value = 0.6855003511296434
$ws.Range("P3").Value = 0.6855003511296432
$ws.Range("Q3").Value = 86.55509297563665
$ws.Range("R3").Value = 519.33055785382
$ws.Range("S3").Value = 0.1051770934051666
$ws.Range("T3").Value = 0.07425000567486074

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.593944666666667
$ws.Range("H4").Value = 10.781834
$ws.Range("I4").Value = 0.02366438041791404
$ws.Range("J4").Value = 0.02505888388006929
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.704179666666667
$ws.Range("N4").Value = 5.112539
$ws.Range("O4").Value = 0.3144996488703566
$ws.Range("P4").Value = 0.3144996488703566
$ws.Range("Q4").Value = 6.124727424058444
$ws.Range("R4").Value = 55.122546816526
$ws.Range("S4").Value = 0.007442439332168509
$ws.Range("T4").Value = 0.007881010181364832

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.593944666666667
$ws.Range("H5").Value = 10.781834
$ws.Range("I5").Value = 0.02366438041791404
$ws.Range("J5").Value = 0.02505888388006929
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.714521666666666
$ws.Range("N5").Value = 11.143565
$ws.Range("O5").Value = 0.6855003511296434
$ws.Range("P5").Value = 0.6855003511296432
$ws.Range("Q5").Value = 13.34978533313444
$ws.Range("R5").Value = 120.14806799821
$ws.Range("S5").Value = 0.01622194108574553
$ws.Range("T5").Value = 0.01717787369870446

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 33.35980533333333
$ws.Range("H6").Value = 100.079416
$ws.Range("I6").Value = 0.2196581186676286
$ws.Range("J6").Value = 0.2326022144589824
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.704179666666667
$ws.Range("N6").Value = 5.112539
$ws.Range("O6").Value = 0.3144996488703566
$ws.Range("P6").Value = 0.3144996488703566
$ws.Range("Q6").Value = 56.85110193302489
$ws.Range("R6").Value = 511.659917397224
$ws.Range("S6").Value = 0.06908240119249233
$ws.Range("T6").Value = 0.07315331477381737

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 33.35980533333333
$ws.Range("H7").Value = 100.079416
$ws.Range("I7").Value = 0.2196581186676286
$ws.Range("J7").Value = 0.2326022144589824
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.714521666666666
$ws.Range("N7").Value = 11.143565
$ws.Range("O7").Value = 0.6855003511296434
$ws.Range("P7").Value = 0.6855003511296432
$ws.Range("Q7").Value = 123.9157197064489
$ws.Range("R7").Value = 1115.24147735804
$ws.Range("S7").Value = 0.1505757174751363
$ws.Range("T7").Value = 0.159448899685165

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 71.394301
$ws.Range("H8").Value = 214.182903
$ws.Range("I8").Value = 0.470096803160314
$ws.Range("J8").Value = 0.4977988434410273
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 1.704179666666667
$ws.Range("N8").Value = 5.112539
$ws.Range("O8").Value = 0.3144996488703566
$ws.Range("P8").Value = 0.3144996488703566
$ws.Range("Q8").Value = 121.6687160800797
$ws.Range("R8").Value = 1095.018444720717
$ws.Range("S8").Value = 0.1478452795289959
$ws.Range("T8").Value = 0.1565575614702727

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 71.394301
$ws.Range("H9").Value = 214.182903
$ws.Range("I9").Value = 0.470096803160314
$ws.Range("J9").Value = 0.4977988434410273
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.714521666666666
$ws.Range("N9").Value = 11.143565
$ws.Range("O9").Value = 0.6855003511296434
$ws.Range("P9").Value = 0.6855003511296432
$ws.Range("Q9").Value = 265.1956779410216
$ws.Range("R9").Value = 2386.761101469195
$ws.Range("S9").Value = 0.3222515236313181
$ws.Range("T9").Value = 0.3412412819707545

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 18.168918
$ws.Range("H10").Value = 54.506754
$ws.Range("I10").Value = 0.1196335022410527
$ws.Range("J10").Value = 0.1266833100162275
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.704179666666667
$ws.Range("N10").Value = 5.112539
$ws.Range("O10").Value = 0.3144996488703566
$ws.Range("P10").Value = 0.3144996488703566
$ws.Range("Q10").Value = 30.963100620934
$ws.Range("R10").Value = 278.667905588406
$ws.Range("S10").Value = 0.03762469444794209
$ws.Range("T10").Value = 0.03984185651783809

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 18.168918
$ws.Range("H11").Value = 54.506754
$ws.Range("I11").Value = 0.1196335022410527
$ws.Range("J11").Value = 0.1266833100162275
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.714521666666666
$ws.Range("N11").Value = 11.143565
$ws.Range("O11").Value = 0.6855003511296434
$ws.Range("P11").Value = 0.6855003511296432
$ws.Range("Q11").Value = 67.48883957088999
$ws.Range("R11").Value = 607.39955613801
$ws.Range("S11").Value = 0.08200880779311057
$ws.Range("T11").Value = 0.08684145349838941

$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 2.052705
$ws.Range("H12").Value = 4.10541
$ws.Range("I12").Value = 0.01351606563570379
$ws.Range("J12").Value = 0.009541696938579768
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.704179666666667
$ws.Range("N12").Value = 5.112539
$ws.Range("O12").Value = 0.3144996488703566
$ws.Range("P12").Value = 0.3144996488703566
$ws.Range("Q12").Value = 3.498178122665
$ws.Range("R12").Value = 20.98906873599
$ws.Range("S12").Value = 0.004250797896537535
$ws.Range("T12").Value = 0.003000860336810694

$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 2.052705
$ws.Range("H13").Value = 4.10541
$ws.Range("I13").Value = 0.01351606563570379
$ws.Range("J13").Value = 0.009541696938579768
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.714521666666666
$ws.Range("N13").Value = 11.143565
$ws.Range("O13").Value = 0.6855003511296434
$ws.Range("P13").Value = 0.6855003511296432
$ws.Range("Q13").Value = 7.624817197774999
$ws.Range("R13").Value = 45.74890318665
$ws.Range("S13").Value = 0.009265267739166253
$ws.Range("T13").Value = 0.006540836601769073
